# Apply updated cryptocurrency price/volume data to Sheet1.
# D-column ("Price") values are forced to Text via a temporary
# "@" NumberFormat (then ClearFormats to avoid leaving style residue)
# so numeric-looking strings (e.g. "1.00", "4.51") are not silently
# coerced into Number cells and lose their original text formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.391.83"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.635.49"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.24%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.57"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.85"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +2.32%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +1.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.03"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +9.30%  "
$ws.Range("E10").Value = "  -1.23%  "
$ws.Range("E11").Value = "  +0.88%  "
$ws.Range("E12").Value = "  +0.34%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.102.90"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +1.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "59.312.11"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "21.27"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +3.37%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000135"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.618.72"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.51"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +3.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "338.65"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.27"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.24"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.31%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.24"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.417"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.95%  "
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.26%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.29"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0751"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  -1.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.87"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "151.49"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.12%  "
$ws.Range("E33").Value = "  +0.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.99"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("E35").Value = "  +2.01%  "
$ws.Range("E36").Value = "  +2.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.837"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("E38").Value = "  -1.03%  "
$ws.Range("E39").Value = "  +1.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "285.11"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.66%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("E44").Value = "  +2.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.06"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0941"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -1.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0227"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.958.47"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.34"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.53"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.85"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.30%  "
